$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.556.88'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '3.124.85'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.118.18'
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('E10').Value = '  +13.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.70'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.469'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.64%  '
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '3.643.03'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.81%  '
$ws.Range('D18').Value = '63.484.70'
$ws.Range('E18').Value = '  +2.60%  '
$ws.Range('D19').Value = '3.125.35'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '463.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.733'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.00%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.110'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.66%  '
$ws.Range('D34').Value = '0.0₃0869'
$ws.Range('E34').Value = '  +8.62%  '
$ws.Range('E35').Value = '  +7.70%  '
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.36'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.09'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '447.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.76'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0372'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '2.893.70'
$ws.Range('E43').Value = '  +1.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.280'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.52%  '
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.06'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.45%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.41%  '
